# Weekly driver report update for 2025-04-21
# Updates the Norma_driver_summary worksheet with refreshed client-count,
# critical-minute, and good-roaming figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers table (rows 3-5) ---
$ws.Range("C3").Value = 91
$ws.Range("D3").Value = 97.5
$ws.Range("C5").Value = 395

# --- Good Drivers table (Total Samples column B) ---
$ws.Range("B15").Value = 56069
$ws.Range("B20").Value = 276086
$ws.Range("B21").Value = 625298
$ws.Range("B26").Value = 331283
$ws.Range("B28").Value = 453652
$ws.Range("B36").Value = 96091
$ws.Range("B39").Value = 99549
$ws.Range("B42").Value = 175767
$ws.Range("B43").Value = 240182
$ws.Range("B51").Value = 684728
$ws.Range("B53").Value = 210188
$ws.Range("B56").Value = 308481
$ws.Range("B63").Value = 443223
$ws.Range("B65").Value = 109665
$ws.Range("B67").Value = 62515
